$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '47.319.04'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '2.490.42'
$ws.Range('E3').Value = '  -0.79%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '321.76'
$ws.Range('E5').Value = '  -0.77%  '
Set-TextValue 'D6' '108.92'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  -0.68%  '
Set-TextValue 'D8' '1.00'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.54%  '
Set-TextValue 'D10' '39.39'
$ws.Range('E10').Value = '  +3.61%  '
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('E12').Value = '  +0.57%  '
Set-TextValue 'D13' '18.53'
$ws.Range('E13').Value = '  +0.57%  '
Set-TextValue 'D14' '7.17'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '2.879.73'
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '2.489.96'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '47.240.36'
$ws.Range('E18').Value = '  -0.95%  '
Set-TextValue 'D19' '13.40'
$ws.Range('E19').Value = '  +5.03%  '
Set-TextValue 'D20' '6.63'
$ws.Range('D21').Value = '0.0₃0939'
$ws.Range('E21').Value = '  +0.03%  '
Set-TextValue 'D22' '2.74'
$ws.Range('E22').Value = '  +14.28%  '
$ws.Range('E23').Value = '  -0.35%  '
Set-TextValue 'D24' '246.86'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  +0.05%  '
Set-TextValue 'D27' '25.68'
$ws.Range('E27').Value = '  -2.53%  '
Set-TextValue 'D28' '2.30'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('E29').Value = '  -0.71%  '
Set-TextValue 'D30' '0.139'
$ws.Range('E30').Value = '  +1.91%  '
Set-TextValue 'D31' '34.68'
$ws.Range('E31').Value = '  -1.71%  '
Set-TextValue 'D32' '49.86'
$ws.Range('E32').Value = '  +0.76%  '
Set-TextValue 'D33' '20.51'
$ws.Range('E33').Value = '  +2.30%  '
Set-TextValue 'D34' '5.31'
$ws.Range('E34').Value = '  -1.23%  '
Set-TextValue 'D35' '0.0785'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  +0.09%  '
Set-TextValue 'D37' '4.78'
$ws.Range('E37').Value = '  +3.08%  '
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E39').Value = '  -2.29%  '
Set-TextValue 'D40' '22.65'
$ws.Range('E40').Value = '  +6.76%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -2.38%  '
Set-TextValue 'D43' '119.38'
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('D45').Value = '1.995.69'
$ws.Range('E45').Value = '  +1.21%  '
Set-TextValue 'D46' '3.02'
$ws.Range('E46').Value = '  +0.36%  '
Set-TextValue 'D47' '2.03'
$ws.Range('E47').Value = '  -3.06%  '
Set-TextValue 'D50' '5.19'
$ws.Range('E50').Value = '  -3.60%  '
Set-TextValue 'D51' '56.79'
$ws.Range('E51').Value = '  +3.16%  '

# Row 48/49: Stacks and FraxShare swap positions with updated price/volume
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D48' '9.09'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D49' '1.78'
$ws.Range('E49').Value = '  -2.16%  '
